$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$data = @(
    @(1, 12345, "3kohm", 150, 1),
    @(2, 12346, "2kohm", 220, 1),
    @(3, 12347, "1.5kohm", 315, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}

$ws.Range("G2").Select()
